# Removed extra datarow and changed to run on Firefox
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SortEmployeeTableColumns")

# Delete the second data row (row 2), shifting all subsequent rows up.
$ws.Rows.Item(2).Delete()

# Update the selection to reflect the new active cell.
$ws.Range("A10").Select()
